$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts old rows 7-10 down to 8-11,
# carries formatting down and auto-adjusts formula ranges, e.g.
# SUM(F2:F7) -> SUM(F2:F8)).
$ws.Rows.Item(7).Insert()

# Fill in the new data row (row 7) with a new time entry.
$ws.Range("A7").Value2 = 2014
$ws.Range("B7").Value2 = 2
$ws.Range("C7").Value2 = 18
$ws.Range("D7").Value2 = 0.625
$ws.Range("E7").Value2 = 0.75
$ws.Range("F7").Formula = "=(E7-D7)*24*60"

# Correct the end time on row 6 (bug fix mentioned in the commit message).
$ws.Range("E6").Value2 = 0.61458333333333337

# Restore the active selection shown in the sheet view.
$ws.Range("A8").Select()
